$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 217; this shifts existing rows 217..315 down to 218..316
# and automatically grows the sheet dimension to A1:R316.
$ws.Rows("217").Insert()

# Populate the newly inserted row 217 with the new data record.
$ws.Range("A217").Value = 3
$ws.Range("B217").Value = "Femacal de La Calera"
$ws.Range("C217").Value = "Coquimbo"
$ws.Range("D217").Value = 44636
$ws.Range("E217").Value = 5
$ws.Range("F217").Value = 100112043
$ws.Range("G217").Value = "Pepino ensalada"
$ws.Range("H217").Value = "Sin especificar"
$ws.Range("I217").Value = "Primera"
$ws.Range("J217").Value = 110
$ws.Range("K217").Value = 16000
$ws.Range("L217").Value = 16500
$ws.Range("M217").Value = 16273
$ws.Range("N217").Value = '$/caja 70 unidades'
$ws.Range("O217").Value = "Provincia de Quillota"
$ws.Range("P217").Value = 232
$ws.Range("Q217").Value = 70
$ws.Range("R217").Value = "Hortaliza"
